$d = $word.ActiveDocument

# Update the date line (single occurrence in the document)
$d.Content.Find.Execute("2023-11-26 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-11-27 Monday", 2) | Out-Null

# Update the division problems in the table. Several expressions (e.g.
# "60÷3=") occur more than once, so each cell is addressed directly by
# row/column and its text is replaced in place rather than relying on a
# document-wide Find/Replace, which would not disambiguate duplicates.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "46÷7="  # was 49÷8=
$t.Cell(1, 2).Range.Text = "36÷2="  # was 38÷7=
$t.Cell(1, 3).Range.Text = "54÷8="  # was 82÷6=
$t.Cell(1, 4).Range.Text = "71÷6="  # was 33÷3=
$t.Cell(1, 5).Range.Text = "15÷5="  # was 77÷8=
$t.Cell(5, 1).Range.Text = "73÷8="  # was 36÷8=
$t.Cell(5, 2).Range.Text = "94÷5="  # was 85÷5=
$t.Cell(5, 3).Range.Text = "38÷6="  # was 58÷9=
$t.Cell(5, 4).Range.Text = "66÷7="  # was 43÷9=
$t.Cell(5, 5).Range.Text = "96÷4="  # was 19÷8=
$t.Cell(9, 1).Range.Text = "50÷4="  # was 60÷3=
$t.Cell(9, 2).Range.Text = "50÷5="  # was 69÷4=
$t.Cell(9, 3).Range.Text = "23÷2="  # was 52÷8=
$t.Cell(9, 4).Range.Text = "77÷9="  # was 60÷3=
$t.Cell(9, 5).Range.Text = "23÷2="  # was 62÷5=
$t.Cell(13, 1).Range.Text = "69÷8="  # was 70÷8=
$t.Cell(13, 2).Range.Text = "57÷8="  # was 33÷8=
$t.Cell(13, 3).Range.Text = "26÷7="  # was 44÷2=
$t.Cell(13, 4).Range.Text = "82÷9="  # was 52÷2=
$t.Cell(13, 5).Range.Text = "21÷6="  # was 28÷5=
$t.Cell(17, 1).Range.Text = "53÷7="  # was 66÷5=
$t.Cell(17, 2).Range.Text = "76÷7="  # was 65÷3=
$t.Cell(17, 3).Range.Text = "96÷8="  # was 76÷9=
$t.Cell(17, 4).Range.Text = "42÷3="  # was 96÷2=
$t.Cell(17, 5).Range.Text = "74÷7="  # was 75÷5=

Write-Host "done"